# Generate Report for Handback
# The "5a0a3879-50bd-4d61-93bc-e59830ad9222" source file has been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet's per-locale status columns, and each locale sheet's
# Status / Latest Target File / Latest Handback File / Latest Handback
# DateTime columns for that row.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$targetFile = "5a0a3879-50bd-4d61-93bc-e59830ad9222.md"

# --- Overview sheet: columns E (zh-cn) and F (de-de) for the 5a0a3879 row (row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack

# --- zh-cn sheet: row 2 is the 5a0a3879 entry
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("I2").Value = $targetFile
$wsZh.Range("J2").Value = $wsZh.Range("G2").Value2
$wsZh.Range("K2").Value = "2016-08-24 20:41:32"

# --- de-de sheet: row 2 is the 5a0a3879 entry
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("I2").Value = $targetFile
$wsDe.Range("J2").Value = $wsDe.Range("G2").Value2
$wsDe.Range("K2").Value = "2016-08-24 20:41:39"

Write-Output "applied handback updates"
